$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 630-631, pushing the existing rows 630-680 down to 632-682.
$ws.Rows("630:631").Insert()

# Populate the newly inserted row 630 with the new weekly data entry
# (Conconina(o), Primera, Región Metropolitana).
$ws.Range("A630").Value = 4
$ws.Range("B630").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C630").Value = "Los Lagos"
$ws.Range("D630").Value = 44826
$ws.Range("E630").Value = 10
$ws.Range("F630").Value = 100112033
$ws.Range("G630").Value = "Lechuga"
$ws.Range("H630").Value = "Conconina(o)"
$ws.Range("I630").Value = "Primera"
$ws.Range("J630").Value = 80
$ws.Range("K630").Value = 14000
$ws.Range("L630").Value = 14000
$ws.Range("M630").Value = 14000
$ws.Range("N630").Value = "$/caja 10 unidades"
$ws.Range("O630").Value = "Región Metropolitana"
$ws.Range("P630").Value = 1400
$ws.Range("Q630").Value = 10
$ws.Range("R630").Value = "Hortaliza"

# Populate the newly inserted row 631 with the new weekly data entry
# (Escarola, Primera, Región de Coquimbo).
$ws.Range("A631").Value = 4
$ws.Range("B631").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C631").Value = "Los Lagos"
$ws.Range("D631").Value = 44826
$ws.Range("E631").Value = 10
$ws.Range("F631").Value = 100112033
$ws.Range("G631").Value = "Lechuga"
$ws.Range("H631").Value = "Escarola"
$ws.Range("I631").Value = "Primera"
$ws.Range("J631").Value = 300
$ws.Range("K631").Value = 15000
$ws.Range("L631").Value = 15000
$ws.Range("M631").Value = 15000
$ws.Range("N631").Value = "$/caja 15 unidades"
$ws.Range("O631").Value = "Región de Coquimbo"
$ws.Range("P631").Value = 1000
$ws.Range("Q631").Value = 15
$ws.Range("R631").Value = "Hortaliza"
